$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '34.438.66'
Set-TextCell 'E2' '  -0.56%  '
Set-TextCell 'D3' '1.803.44'
Set-TextCell 'E3' '  +0.41%  '
Set-TextCell 'E4' '  +0.15%  '
Set-TextCell 'D5' '228.29'
Set-TextCell 'E5' '  +0.55%  '
Set-TextCell 'E6' '  +4.40%  '
Set-TextCell 'E7' '  +0.18%  '
Set-TextCell 'D8' '34.77'
Set-TextCell 'E8' '  +5.27%  '
Set-TextCell 'E9' '  +1.23%  '
Set-TextCell 'D10' '0.0694'
Set-TextCell 'E10' '  -0.18%  '
Set-TextCell 'E11' '  +0.16%  '
Set-TextCell 'D12' '2.065.48'
Set-TextCell 'E12' '  +0.54%  '
Set-TextCell 'D13' '11.21'
Set-TextCell 'E13' '  +0.90%  '
Set-TextCell 'D14' '1.800.00'
Set-TextCell 'E14' '  +0.21%  '
Set-TextCell 'D15' '0.643'
Set-TextCell 'E15' '  +0.87%  '
Set-TextCell 'D16' '34.442.84'
Set-TextCell 'E16' '  -0.35%  '
Set-TextCell 'D17' '4.37'
Set-TextCell 'E17' '  +1.78%  '
Set-TextCell 'D18' '69.08'
Set-TextCell 'E18' '  +0.26%  '
Set-TextCell 'D19' '0.0₃0798'
Set-TextCell 'E19' '  -0.53%  '
Set-TextCell 'D20' '245.68'
Set-TextCell 'E20' '  -0.91%  '
Set-TextCell 'D21' '11.50'
Set-TextCell 'E21' '  +1.69%  '
Set-TextCell 'E22' '  +0.11%  '
Set-TextCell 'D23' '4.16'
Set-TextCell 'E23' '  -0.57%  '
Set-TextCell 'D24' '173.11'
Set-TextCell 'E24' '  +3.87%  '
Set-TextCell 'E25' '  +1.43%  '
Set-TextCell 'D26' '7.85'
Set-TextCell 'E26' '  +7.36%  '
Set-TextCell 'B27' 'EthereumClassic'
Set-TextCell 'C27' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 'D27' '16.79'
Set-TextCell 'E27' '  +1.19%  '
Set-TextCell 'B28' 'Stellar'
Set-TextCell 'C28' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D28' '0.119'
Set-TextCell 'E28' '  +2.46%  '
Set-TextCell 'E29' '  -0.07%  '
Set-TextCell 'D30' '4.00'
Set-TextCell 'E30' '  -2.62%  '
Set-TextCell 'D31' '0.0530'
Set-TextCell 'E31' '  +0.76%  '
Set-TextCell 'D32' '3.83'
Set-TextCell 'E32' '  +0.39%  '
Set-TextCell 'D33' '1.24'
Set-TextCell 'E33' '  +0.20%  '
Set-TextCell 'E34' '  -0.02%  '
Set-TextCell 'D35' '1.394.58'
Set-TextCell 'E35' '  -2.41%  '
Set-TextCell 'D36' '0.680'
Set-TextCell 'E36' '  +1.15%  '
Set-TextCell 'D37' '2.51'
Set-TextCell 'E37' '  -3.92%  '
Set-TextCell 'D38' '1.06'
Set-TextCell 'E38' '  -0.50%  '
Set-TextCell 'E39' '  -1.56%  '
Set-TextCell 'D40' '83.38'
Set-TextCell 'E40' '  -2.73%  '
Set-TextCell 'E41' '  +2.68%  '
Set-TextCell 'D42' '0.949'
Set-TextCell 'E42' '  +1.48%  '
Set-TextCell 'E43' '  -0.44%  '
Set-TextCell 'D44' '13.53'
Set-TextCell 'E44' '  -0.65%  '
Set-TextCell 'E45' '  +3.69%  '
Set-TextCell 'E46' '  -3.12%  '
Set-TextCell 'E47' '  -2.40%  '
Set-TextCell 'D48' '1.964.61'
Set-TextCell 'E48' '  +0.55%  '
Set-TextCell 'D49' '104.93'
Set-TextCell 'E49' '  -1.30%  '
Set-TextCell 'E50' '  +0.16%  '
Set-TextCell 'D51' '0.0₆0129'
Set-TextCell 'E51' '  +0.63%  '
